$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # "Proveedores"
$ws2 = $wb.Worksheets.Item(2)  # "Productos"

# The block of data that used to live in H13:J18 (separated from the main
# table by a gap) is moved up to sit right next to the main table, in H1:J6.
$src = $ws1.Range("H13:J18")
$dst = $ws1.Range("H1")
$src.Cut($dst)

# The rows that used to hold that block (and the blank spacer row 12) are
# now completely empty, so remove them outright.
$ws1.Rows("7:18").Delete()

# Update the remembered selection on each sheet and make "Proveedores" the
# active (selected) tab instead of "Productos".
$ws1.Range("J1:J6").Select()
$ws2.Range("J1:J4").Select()
$ws1.Activate()
